# "10Th - MB for single stock and added new group"
#
# The weekly MarketBeat rank sheet gets a new week's worth of columns
# inserted at the front (pushing the previous weeks right), the new
# week headers get filled in, the new week's analyst cells are filled
# with the usual "UN" placeholder, and two new research firms
# (Benchmark, Evercore ISI) are appended as new rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns in front of column B. This shifts the existing
# B:E (Jun_17, Jun_15, Jun_13, Jun_10) right to E:H, carrying each
# column's cell styles (including the highlighted "latest action"
# cells) along with it.
$ws.Range("B1:D1").EntireColumn.Insert()

# New week headers in row 1.
$ws.Cells.Item(1, 2).Value = "Jun_27"
$ws.Cells.Item(1, 3).Value = "Jun_26"
$ws.Cells.Item(1, 4).Value = "Jun_26"

# Fill the new B/C/D columns for every existing analyst row with the
# usual "UN" (no rating change) placeholder used throughout the sheet.
for ($r = 2; $r -le 27; $r++) {
  $ws.Cells.Item($r, 2).Value = "UN"
  $ws.Cells.Item($r, 3).Value = "UN"
  $ws.Cells.Item($r, 4).Value = "UN"
}

# New research-firm group added at the bottom of the sheet.
$ws.Cells.Item(28, 1).Value = "Benchmark"
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"
$ws.Cells.Item(28, 4).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI"
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
$ws.Cells.Item(29, 4).Value = "UN"
